$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the DLC table with a new row (11 m/s turbulent case), cloned from
# row 4 so formatting/styles carry over, then adjust formulas/values that
# differ for the new case.
$ws.Range("A4:AM4").Copy($ws.Range("A5"))

$ws.Range("A5").Formula = "=A4"
$ws.Range("B5").Formula = "=A5+20"
$ws.Range("D5").Formula = "=""dlc01_steady_wsp"" & E5 & ""_s101"""
$ws.Range("E5").Formula = "=E4+1"
$ws.Range("I5").Formula = "=(0.16*(0.75*E5+5.6))/E5"
$ws.Range("J5").Value = "turb_s101_11ms"
$ws.Range("K5").Formula = "=E5*B5/512"
$ws.Range("M5").Formula = "=8/E5"

$ws.Range("A5").Select()

# Cosmetic: the tab-split ratio was also nudged in the authored workbook.
$excel.ActiveWindow.TabRatio = 0.5
